$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 00:53"

# Reorder countries (swap text values) to reflect updated ranking
$ws.Range("A10").Value = "Colombia"
$ws.Range("A11").Value = "Mexico"
$ws.Range("A124").Value = "Tunez"
$ws.Range("A125").Value = "Mayotte"
$ws.Range("A126").Value = "Somalia"
$ws.Range("A151").Value = "Trinidad yTobago"
$ws.Range("A152").Value = "Georgia"
$ws.Range("A159").Value = "Guyana"
$ws.Range("A160").Value = "Principado de Andorra"

# Update numeric statistics
$ws.Range("B4").Value = 6039973
$ws.Range("C4").Value = 39608
$ws.Range("D4").Value = 3339932
$ws.Range("E4").Value = 2515389
$ws.Range("G4").Value = 999
$ws.Range("H4").Value = 184652

$ws.Range("D5").Value = 2947250
$ws.Range("E5").Value = 695492

$ws.Range("B10").Value = 582022
$ws.Range("C10").Value = 9752
$ws.Range("D10").Value = 417793
$ws.Range("E10").Value = 145761
$ws.Range("G10").Value = 284
$ws.Range("H10").Value = 18468

$ws.Range("B11").Value = 573888
$ws.Range("C11").Value = 5267
$ws.Range("D11").Value = 396758
$ws.Range("E11").Value = 115054
$ws.Range("G11").Value = 626
$ws.Range("H11").Value = 62076

$ws.Range("B27").Value = 126780
$ws.Range("C27").Value = 363
$ws.Range("D27").Value = 112771
$ws.Range("E27").Value = 4908
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 9101

$ws.Range("B34").Value = 98062
$ws.Range("C34").Value = 237
$ws.Range("D34").Value = 69612
$ws.Range("E34").Value = 23108
$ws.Range("G34").Value = 25
$ws.Range("H34").Value = 5342

$ws.Range("B48").Value = 64668
$ws.Range("C48").Value = 846
$ws.Range("D48").Value = 52823
$ws.Range("E48").Value = 10619
$ws.Range("G48").Value = 17
$ws.Range("H48").Value = 1226

$ws.Range("B53").Value = 53317
$ws.Range("C53").Value = 296
$ws.Range("D53").Value = 40726
$ws.Range("E53").Value = 11580
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 1011

$ws.Range("B76").Value = 19142
$ws.Range("C76").Value = 169
$ws.Range("D76").Value = 17651
$ws.Range("E76").Value = 1080
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 411

$ws.Range("B81").Value = 15908
$ws.Range("C81").Value = 157
$ws.Range("D81").Value = 11044
$ws.Range("E81").Value = 4270
$ws.Range("G81").Value = 8
$ws.Range("H81").Value = 594

$ws.Range("B86").Value = 13294
$ws.Range("C86").Value = 108
$ws.Range("D86").Value = 8974
$ws.Range("E86").Value = 4043
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 277

$ws.Range("B90").Value = 10542
$ws.Range("C90").Value = 38
$ws.Range("E90").Value = 930

$ws.Range("B100").Value = 8151
$ws.Range("C100").Value = 29
$ws.Range("D100").Value = 5743
$ws.Range("E100").Value = 2208
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 200

$ws.Range("B124").Value = 3323
$ws.Range("C124").Value = 117
$ws.Range("D124").Value = 1504
$ws.Range("E124").Value = 1746
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 73

$ws.Range("B125").Value = 3301
$ws.Range("C125").Value = 64
$ws.Range("D125").Value = 2964
$ws.Range("E125").Value = 297
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 40

$ws.Range("B126").Value = 3275
$ws.Range("D126").Value = 2443
$ws.Range("E126").Value = 737
$ws.Range("H126").Value = 95

$ws.Range("B151").Value = 1476
$ws.Range("C151").Value = 65
$ws.Range("D151").Value = 594
$ws.Range("E151").Value = 867
$ws.Range("H151").Value = 15

$ws.Range("B152").Value = 1447
$ws.Range("C152").Value = 11
$ws.Range("D152").Value = 1190
$ws.Range("E152").Value = 238
$ws.Range("H152").Value = 19

$ws.Range("B159").Value = 1140
$ws.Range("C159").Value = 47
$ws.Range("D159").Value = 616
$ws.Range("E159").Value = 492
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 32

$ws.Range("B160").Value = 1098
$ws.Range("D160").Value = 893
$ws.Range("E160").Value = 152
$ws.Range("H160").Value = 53

$ws.Range("B165").Value = 894
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 846
$ws.Range("E165").Value = 33

$ws.Range("B180").Value = 355
$ws.Range("C180").Value = 1
$ws.Range("E180").Value = 10

